# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker data table (rows 16-25, columns B:G) is re-sorted: instead of
# being grouped by worker (all periods for ANGIE, then all periods for
# ANIBAL), it is now grouped by period first (1711, 1712, 1801, 1802, 1803),
# with each worker's row for that period following.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CC", "1143387210", "ANGIE STEPHANIE CAMACHO AGRESOTT", "1711", 29509, 737717),
    @("CC", "1047423122", "ANIBAL FELIPE ARENAS SOTO",        "1711", 29509, 737717),
    @("CC", "1143387210", "ANGIE STEPHANIE CAMACHO AGRESOTT", "1712", 29509, 737717),
    @("CC", "1047423122", "ANIBAL FELIPE ARENAS SOTO",        "1712", 29509, 737717),
    @("CC", "1143387210", "ANGIE STEPHANIE CAMACHO AGRESOTT", "1801", 29509, 737717),
    @("CC", "1047423122", "ANIBAL FELIPE ARENAS SOTO",        "1801", 29509, 737717),
    @("CC", "1143387210", "ANGIE STEPHANIE CAMACHO AGRESOTT", "1802", 29509, 737717),
    @("CC", "1047423122", "ANIBAL FELIPE ARENAS SOTO",        "1802", 29509, 737717),
    @("CC", "1143387210", "ANGIE STEPHANIE CAMACHO AGRESOTT", "1803", 19673, 737717),
    @("CC", "1047423122", "ANIBAL FELIPE ARENAS SOTO",        "1803", 29509, 737717)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
}
